$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "47.388.51"
Set-TextValue "E2" "  +2.99%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.503.95"
Set-TextValue "E3" "  +2.44%  "

# Row 5 - BNB
Set-TextValue "D5" "324.87"
Set-TextValue "E5" "  +1.25%  "

# Row 6 - Solana
Set-TextValue "D6" "109.99"
Set-TextValue "E6" "  +5.25%  "

# Row 7 - XRP
Set-TextValue "E7" "  +1.32%  "

# Row 8 - USDC
Set-TextValue "E8" "  +0.00%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +0.44%  "

# Row 10 - Avalanche
Set-TextValue "D10" "39.56"
Set-TextValue "E10" "  +10.26%  "

# Row 11 - Dogecoin
Set-TextValue "E11" "  +1.32%  "

# Row 12 - TRON
Set-TextValue "E12" "  +0.88%  "

# Row 13 - Chainlink
Set-TextValue "D13" "18.54"
Set-TextValue "E13" "  +1.93%  "

# Row 14 - Polkadot
Set-TextValue "E14" "  +2.39%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.895.51"
Set-TextValue "E15" "  +2.45%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "2.510.32"
Set-TextValue "E16" "  +2.51%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.863"
Set-TextValue "E17" "  +2.67%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "47.334.20"
Set-TextValue "E18" "  +3.25%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue "D19" "12.88"
Set-TextValue "E19" "  +2.98%  "

# Row 20 - Uniswap
Set-TextValue "E20" "  +4.12%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0943"
Set-TextValue "E21" "  +1.09%  "

# Row 22 - ImmutableX
Set-TextValue "E22" "  +13.55%  "

# Row 23 - Litecoin
Set-TextValue "D23" "70.56"
Set-TextValue "E23" "  -0.90%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "248.80"
Set-TextValue "E24" "  +0.90%  "

# Row 25 - PancakeSwap
Set-TextValue "E25" "  +3.77%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "26.12"
Set-TextValue "E26" "  +1.06%  "

# Row 27 - Dai
Set-TextValue "E27" "  -0.04%  "

# Row 28 - was Cosmos, now Toncoin
Set-TextValue "B28" "Toncoin"
Set-TextValue "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "2.30"
Set-TextValue "E28" "  +1.85%  "

# Row 29 - was Toncoin, now Cosmos
Set-TextValue "B29" "Cosmos"
Set-TextValue "C29" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "10.08"
Set-TextValue "E29" "  +4.29%  "

# Row 30 - InjectiveProtocol
Set-TextValue "D30" "35.45"
Set-TextValue "E30" "  +5.19%  "

# Row 31 - Kaspa
Set-TextValue "E31" "  +9.65%  "

# Row 32 - OKB
Set-TextValue "D32" "49.93"
Set-TextValue "E32" "  +1.38%  "

# Row 33 - Celestia
Set-TextValue "D33" "19.99"
Set-TextValue "E33" "  -0.21%  "

# Row 34 - Filecoin
Set-TextValue "E34" "  +1.57%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.0797"
Set-TextValue "E35" "  +5.03%  "

# Row 36 - FirstDigitalUSD
Set-TextValue "E36" "  +0.32%  "

# Row 37 - ARBITRUM
Set-TextValue "E37" "  +6.14%  "

# Row 38 - RenderToken
Set-TextValue "E38" "  +4.35%  "

# Row 39 - LidoDAOToken
Set-TextValue "E39" "  +1.95%  "

# Row 40 - Stellar
Set-TextValue "E40" "  +1.73%  "

# Row 41 - Monero
Set-TextValue "D41" "121.62"
Set-TextValue "E41" "  -3.73%  "

# Row 42 - WEMIXToken
Set-TextValue "E42" "  -1.36%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "21.26"
Set-TextValue "E43" "  +1.57%  "

# Row 44 - VeChain
Set-TextValue "D44" "0.0300"
Set-TextValue "E44" "  +2.49%  "

# Row 45 - Maker
Set-TextValue "D45" "2.001.27"
Set-TextValue "E45" "  +2.18%  "

# Row 46 - NEARProtocol
Set-TextValue "D46" "3.12"
Set-TextValue "E46" "  +5.49%  "

# Row 47 - ApeXProtocol
Set-TextValue "E47" "  -1.97%  "

# Row 48 - Stacks
Set-TextValue "E48" "  -3.47%  "

# Row 49 - FraxShare
Set-TextValue "E49" "  -0.43%  "

# Row 50 - THORChain
Set-TextValue "D50" "5.24"
Set-TextValue "E50" "  +4.95%  "

# Row 51 - MultiversX
Set-TextValue "D51" "56.82"
Set-TextValue "E51" "  +4.45%  "
